# "table removed from script template"
#
# The single paragraph holding the "$characters" placeholder (styled
# "Character") is replaced by three paragraphs, each styled "Character"
# and laid out as two tab-separated placeholders:
#   $char00 <tab> $char01
#   $char10 <tab> $char11
#   $char20 <tab> $char21
# The trailing _GoBack bookmark that used to sit at the end of the
# original paragraph is preserved at the end of the new, final paragraph.

$d = $word.ActiveDocument

# Locate the paragraph that still contains the "$characters" placeholder.
# Find.Execute collapses $searchRange to the matched text in place, so we
# must keep using that same range object (not re-derive a fresh one from
# $d.Content, which would span the whole document again).
$searchRange = $d.Content
$found = $searchRange.Find.Execute("`$characters", $false, $false, $false,
                                    $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the `$characters placeholder paragraph"
}

$targetParagraph = $searchRange.Paragraphs(1)
$targetRange = $targetParagraph.Range

$wordMlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$newParagraphsXml =
    '<w:p ' + $wordMlNs + '><w:pPr><w:pStyle w:val="Character"/></w:pPr>' +
    '<w:r><w:t>$char00</w:t></w:r><w:r><w:tab/></w:r><w:r><w:t>$char01</w:t></w:r></w:p>' +
    '<w:p ' + $wordMlNs + '><w:pPr><w:pStyle w:val="Character"/></w:pPr>' +
    '<w:r><w:t>$char10</w:t></w:r><w:r><w:tab/></w:r><w:r><w:t>$char11</w:t></w:r></w:p>' +
    '<w:p ' + $wordMlNs + '><w:pPr><w:pStyle w:val="Character"/></w:pPr>' +
    '<w:r><w:t>$char20</w:t></w:r><w:r><w:tab/></w:r><w:r><w:t>$char21</w:t></w:r>' +
    '<w:bookmarkStart w:id="1" w:name="_GoBack"/><w:bookmarkEnd w:id="1"/></w:p>'

$targetRange.InsertXML($newParagraphsXml)
